# Start of 2025 author list / membership refresh: add Frederic Pouzoulet as a new
# collaborator row, and update Alfredo Fernandez-Rodriguez's contact details so the
# JISCMail and collaboration database entries are consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Alfredo Fernandez-Rodriguez) ---
# "Name on Publications" is filled in, and the email address gains a second address.
$ws.Range("F2").Value = "Alfredo Fernandez-Rodriguez"
$ws.Range("E2").Value = "Alfredo.fernandez-rodriguez@curie.fr; alfredofernandezrod@gmail.com"

# --- Add new row 3 (Frederic Pouzoulet) ---
$ws.Range("E3").Value = "frederic.pouzoulet@curie.fr"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:frederic.pouzoulet@curie.fr") | Out-Null

$ws.Range("A3").Value = "Prof."
$ws.Range("B3").Value = "Frederic"
$ws.Range("C3").Value = "Pouzoulet"
$ws.Range("D3").Value = "F."
$ws.Range("F3").Value = "F.Pouzoulet"
$ws.Range("G3").Value = "Inst-Curie"
$ws.Range("H3").Value = "Institut Curie-Orsay Research Center, Bat a Campus d'Orsay, 91400 Orsay, France"
$ws.Range("I3").Value = 0

# Grow the table so the new row is included
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:P3"))

# Match the saved cursor position
$ws.Range("J3").Select() | Out-Null
